# Adds a new "2022-Q3" quarter to the 00998-中信银行 workbook:
#  1. Insert a new sheet "2022-Q3" right after "总计" and before "2022-Q2",
#     populated with the fund-holding detail table for the new quarter.
#  2. Insert a new row at the top of the "总计" summary table with the
#     2022-Q3 aggregate figures, pushing the existing rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Build the new "2022-Q3" worksheet
# ---------------------------------------------------------------------------
$summarySheet = $wb.Worksheets.Item(1)      # "总计"
$q2Sheet      = $wb.Worksheets.Item(2)      # "2022-Q2" (will stay right after the new sheet)

$newSheet = $wb.Worksheets.Add($q2Sheet)    # inserted before "2022-Q2"
$newSheet.Name = "2022-Q3"

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Fund holding rows (A=index, B=code, C=name, D..G=text-like numeric fields, H=rank number)
$rows = @(
    @(0,  "008515", "富兰克林国海基本面优选混合",                 "11.82", "85.30", "7.20", "0.8510", 1),
    @(1,  "011152", "富兰克林国海兴海回报混合",                   "14.87", "84.99", "4.55", "0.6766", 6),
    @(2,  "513690", "博时恒生港股通高股息率ETF",                  "3.05",  "97.26", "2.64", "0.0805", 5),
    @(3,  "501305", "汇添富中证港股通高股息投资指数（LOF）A",       "0.87",  "92.21", "3.41", "0.0297", 7),
    @(4,  "513530", "华泰柏瑞中证港股通高股息投资ETF（QDII）",      "0.78",  "95.80", "3.65", "0.0285", 7),
    @(5,  "159726", "华夏恒生中国内地企业高股息率ETF",             "0.84",  "96.48", "2.91", "0.0244", 4),
    @(6,  "004098", "前海开源港股通股息率50强股票",               "0.31",  "87.14", "3.72", "0.0115", 3),
    @(7,  "007751", "景顺长城中证沪港深红利成长低波动指数A",        "0.67",  "90.27", "1.66", "0.0111", 6),
    @(8,  "501306", "汇添富中证港股通高股息投资指数（LOF）C",       "0.17",  "92.21", "3.41", "0.0058", 7),
    @(9,  "006106", "景顺长城量化港股通股票",                     "0.52",  "50.20", "1.11", "0.0058", 8),
    @(10, "005702", "恒生前海港股通高股息低波动指数",              "0.20",  "94.22", "2.57", "0.0051", 4),
    @(11, "006658", "财通中证香港红利等权投资指数A",               "0.13",  "88.79", "3.55", "0.0046", 9),
    @(12, "006659", "财通中证香港红利等权投资指数C",               "0.04",  "88.79", "3.55", "0.0014", 9),
    @(13, "007760", "景顺长城中证沪港深红利成长低波动指数C",        "0.06",  "90.27", "1.66", "0.0010", 6)
)

# Force columns B, D, E, F and G to be stored as text so that fund codes keep
# their leading zeros and the decimal values keep their exact printed form
# (e.g. "85.30", "0.8510") instead of being reinterpreted as numbers.
$lastRow = 1 + $rows.Count
$newSheet.Range("B2:B$lastRow").NumberFormat = "@"
$newSheet.Range("D2:G$lastRow").NumberFormat = "@"

$r = 2
foreach ($row in $rows) {
    $newSheet.Range("A$r").Value = $row[0]
    $newSheet.Range("B$r").Value = $row[1]
    $newSheet.Range("C$r").Value = $row[2]
    $newSheet.Range("D$r").Value = $row[3]
    $newSheet.Range("E$r").Value = $row[4]
    $newSheet.Range("F$r").Value = $row[5]
    $newSheet.Range("G$r").Value = $row[6]
    $newSheet.Range("H$r").Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2. Insert the 2022-Q3 row into the "总计" summary sheet
# ---------------------------------------------------------------------------
$summarySheet.Rows.Item(2).Insert()

# Clear the formatting the Insert operation copied down from the header row
# for columns B..D (they should stay unstyled, like all the other data rows).
$summarySheet.Range("B2:D2").ClearFormats()

# Re-apply column A's usual style (bold / bordered) by copying it from the
# row that used to be the first data row and is now row 3.
$summarySheet.Range("A3").Copy()
$summarySheet.Range("A2").PasteSpecial(-4122)   # xlPasteFormats

$summarySheet.Range("A2").Value = 0
$summarySheet.Range("B2").Value = "2022-Q3"
$summarySheet.Range("C2").Value = 14
$summarySheet.Range("D2").Value = 1.74

# Leave the workbook's active sheet as "总计", matching the original file.
$summarySheet.Activate()
